$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp refreshed from 13:42 to 14:12
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 14:12"

# Provincia case-count table refreshed with newer figures; rows stay sorted
# by "Casos totales" (column B) descending, so many province labels shift
# down/up a row or two along with their B:E statistics.
$data = @(
    @(9, "Navarra", 1411, 35, 1327, 49),
    @(10, "Araba/Alava", 1207, 466, 1126, 81),
    @(11, "Ciudad Real", 1147, 95, 1050, 89),
    @(12, "Alacant/Alicante", 1039, 17, 941, 81),
    @(13, "La Rioja", 995, 48, 904, 43),
    @(14, "Toledo", 965, 95, 860, 78),
    @(15, "Aragon", 907, 29, 838, 40),
    @(16, "Malaga", 819, 48, 731, 40),
    @(17, "A Coruña", 812, 47, 789, 23),
    @(18, "Asturias", 779, 35, 719, 25),
    @(19, "Pontevedra", 689, 47, 684, 5),
    @(20, "Cantabria", 671, 14, 640, 17),
    @(21, "Albacete", 666, 95, 592, 66),
    @(22, "Salamanca", 629, 73, 502, 54),
    @(23, "Valladolid", 598, 50, 521, 27),
    @(24, "Murcia", 596, 9, 579, 8),
    @(25, "Granada", 579, 1, 553, 25),
    @(26, "Gipuzkoa/Guipuzcoa", 563, 466, 543, 20),
    @(27, "Tenerife", 539, 15, 519, 24),
    @(28, "Sevilla", 535, 6, 511, 18),
    @(29, "Burgos", 530, 67, 431, 32),
    @(30, "Leon", 516, 48, 422, 46),
    @(31, "Caceres", 485, 3, 447, 35),
    @(32, "Guadalajara", 428, 95, 370, 56),
    @(34, "Segovia", 361, 62, 262, 37),
    @(35, "Soria", 339, 32, 291, 16),
    @(36, "Zaragoza", 329, 0, 315, 14),
    @(37, "Jaen", 316, 5, 297, 14),
    @(38, "Cordoba", 291, 0, 285, 6),
    @(39, "Cadiz", 278, 4, 270, 4),
    @(40, "Avila", 270, 53, 190, 27),
    @(41, "Badajoz", 257, 5, 248, 4),
    @(42, "Ourense", 235, 47, 231, 4),
    @(44, "Cuenca", 177, 95, 142, 27),
    @(46, "Palencia", 139, 14, 120, 5),
    @(47, "Lugo", 132, 47, 128, 4),
    @(48, "Almeria", 115, 5, 105, 5),
    @(49, "Zamora", 106, 24, 73, 9)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
